# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# This updates the "K" column (column G) values for rows 2-46 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G (K), rows 2 through 46, in order.
$newValues = @(1, 0, 0, 2, 1, 1, 3, 0, 1, 2, 0, 1, 0, 2, 4, 4, 1, 3, 0, 2, 2, 1, 0, 0, 1, 3, 0, 2, 2, 4, 2, 2, 1, 1, 1, 2, 1, 1, 2, 0, 0, 2, 2, 1, 0)

$startRow = 2
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $newValues[$i]
}
